# ---------------------------------------------------------------------------
# Commit: "Add upgrading UI and tower attributes UI"
#
# 1) Refresh the cached "datetimeFigureOut" field text (10/27/2023 -> 11/2/2023)
#    on the slide master and every slide layout (PowerPoint re-caches this
#    field's rendered text whenever the deck is re-saved on a later day).
# 2) Slide 1: nudge/resize the four tower-icon pictures (תמונה 13-16).
# 3) Slide 2: reflow the "tower attributes" cluster - reposition/resize the
#    four icon pictures (תמונה 1, 2, 3, 6), crop three of them, rename +
#    reposition the rounded-rectangle badge behind them, and re-stack it in
#    front of the icons (it used to be drawn first/behind; now it is drawn
#    last/in front).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Date placeholder text on master + every layout ---------------------

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "10/27/2023") {
                $tr.Text = "11/2/2023"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DateShapes $layout.Shapes
}

# --- 2) Slide 1: tower icon pictures ---------------------------------------

$slide1 = $p.Slides.Item(1)

$t13 = $slide1.Shapes.Item("תמונה 13")
$t13.Left   = 121.86205291748047
$t13.Top    = 268.24725341796875
$t13.Width  = 24.09449005126953
$t13.Height = 30.047245025634766

$t14 = $slide1.Shapes.Item("תמונה 14")
$t14.Left   = 166.6167755126953
$t14.Top    = 268.24749755859375
$t14.Width  = 24.09449005126953
$t14.Height = 30.047245025634766

$t15 = $slide1.Shapes.Item("תמונה 15")
$t15.Left   = 212.06370544433594
$t15.Top    = 268.6013488769531
$t15.Width  = 24.09449005126953
$t15.Height = 30.047245025634766

$t16 = $slide1.Shapes.Item("תמונה 16")
$t16.Left   = 257.0809631347656
$t16.Top    = 268.24749755859375
$t16.Width  = 24.09449005126953
$t16.Height = 30.026378631591797

# --- 3) Slide 2: tower attribute icons + badge ------------------------------

$slide2 = $p.Slides.Item(2)

$icon1 = $slide2.Shapes.Item("תמונה 1")
$icon1.Left   = 230.99119567871094
$icon1.Top    = 272.4375
$icon1.Width  = 21.838977813720703
$icon1.Height = 26.567401885986328
$icon1.PictureFormat.CropTop = 0
$icon1.PictureFormat.CropBottom = 0

$icon2 = $slide2.Shapes.Item("תמונה 2")
$icon2.Left   = 269.5804748535156
$icon2.Top    = 272.6249694824219
$icon2.Width  = 21.838977813720703
$icon2.Height = 26.40858268737793
$icon2.PictureFormat.CropTop = 0.36855

$icon3 = $slide2.Shapes.Item("תמונה 3")
$icon3.Left   = 308.1697692871094
$icon3.Top    = 267.75
$icon3.Width  = 22.066064834594727
$icon3.Height = 31.2548828125
$icon3.PictureFormat.CropTop = 0.1249875

$icon6 = $slide2.Shapes.Item("תמונה 6")
$icon6.Left   = 346.98614501953125
$icon6.Top    = 264.4666442871094
$icon6.Width  = 22.066064834594727
$icon6.Height = 34.53826904296875

$badge = $slide2.Shapes.Item("מלבן מעוגל 4")
$badge.Name   = "מלבן מעוגל 18"
$badge.Left   = 193.87496948242188
$badge.Top    = 276.75
$badge.Width  = 21.727638244628906
$badge.Height = 21.826772689819336
# Bring the badge to the front so it is stacked after the icons (it used to
# be drawn first/behind them).
$badge.ZOrder(0)
